$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 327307
$ws.Range("D2").Value = 416903939
$ws.Range("C3").Value = 264
$ws.Range("D3").Value = 316083
$ws.Range("C4").Value = 329
$ws.Range("D4").Value = 470692
$ws.Range("C8").Value = 875
$ws.Range("D8").Value = 1286899
$ws.Range("C10").Value = 118799
$ws.Range("D10").Value = 174059813
$ws.Range("C12").Value = 60677
$ws.Range("D12").Value = 87562990
$ws.Range("C14").Value = 50
$ws.Range("D14").Value = 69039
$ws.Range("C16").Value = 4046
$ws.Range("D16").Value = 5743697
$ws.Range("C19").Value = 74
$ws.Range("D19").Value = 107535
$ws.Range("C20").Value = 7034
$ws.Range("D20").Value = 9819014
$ws.Range("C22").Value = 78799
$ws.Range("D22").Value = 98176338
$ws.Range("C28").Value = 32862
$ws.Range("D28").Value = 48096936
$ws.Range("C30").Value = 11676
$ws.Range("D30").Value = 16796116
$ws.Range("C33").Value = 1572
$ws.Range("D33").Value = 2208781
$ws.Range("C35").Value = 1935
$ws.Range("D35").Value = 2731421
$ws.Range("C36").Value = 98731
$ws.Range("D36").Value = 124153656
$ws.Range("C42").Value = 906
$ws.Range("D42").Value = 1333685
$ws.Range("C44").Value = 44882
$ws.Range("D44").Value = 65771629
$ws.Range("C46").Value = 9314
$ws.Range("D46").Value = 13357668
$ws.Range("C48").Value = 1425
$ws.Range("D48").Value = 1979811
$ws.Range("C51").Value = 2462
$ws.Range("D51").Value = 3441802
$ws.Range("C52").Value = 70212
$ws.Range("D52").Value = 88044381
$ws.Range("C59").Value = 28575
$ws.Range("D59").Value = 41906493
$ws.Range("C62").Value = 11351
$ws.Range("D62").Value = 16410368
$ws.Range("C64").Value = 1371
$ws.Range("D64").Value = 1916737
$ws.Range("C68").Value = 1562
$ws.Range("D68").Value = 2188581
$ws.Range("C70").Value = 20840
$ws.Range("D70").Value = 27287296
$ws.Range("C74").Value = 7705
$ws.Range("D74").Value = 11284101
$ws.Range("C76").Value = 5208
$ws.Range("D76").Value = 7562986
$ws.Range("C77").Value = 496
$ws.Range("D77").Value = 702739
$ws.Range("C79").Value = 143405
$ws.Range("D79").Value = 178706970
$ws.Range("C83").Value = 435
$ws.Range("D83").Value = 635324
$ws.Range("C85").Value = 64437
$ws.Range("D85").Value = 94431439
$ws.Range("C88").Value = 30268
$ws.Range("D88").Value = 43785165
$ws.Range("C90").Value = 2767
$ws.Range("D90").Value = 3984198
$ws.Range("C91").Value = 2971
$ws.Range("D91").Value = 4200764
$ws.Range("C92").Value = 34175
$ws.Range("D92").Value = 46334405
$ws.Range("C96").Value = 8270
$ws.Range("D96").Value = 12156583
$ws.Range("C98").Value = 7639
$ws.Range("D98").Value = 11083971
$ws.Range("C100").Value = 545
$ws.Range("D100").Value = 773406
$ws.Range("C101").Value = 514
$ws.Range("D101").Value = 742026
$ws.Range("C102").Value = 11079
$ws.Range("D102").Value = 17392189
$ws.Range("C104").Value = 2716
$ws.Range("D104").Value = 4577049
$ws.Range("C106").Value = 3695
$ws.Range("D106").Value = 6211094
$ws.Range("C108").Value = 163
$ws.Range("D108").Value = 271445
$ws.Range("C109").Value = 213
$ws.Range("D109").Value = 339180
$ws.Range("C110").Value = 144204
$ws.Range("D110").Value = 178353166
$ws.Range("C114").Value = 960
$ws.Range("D114").Value = 1408477
$ws.Range("C116").Value = 53457
$ws.Range("D116").Value = 78347067
$ws.Range("C118").Value = 27749
$ws.Range("D118").Value = 40203803
$ws.Range("C119").Value = 1323
$ws.Range("D119").Value = 1809431
$ws.Range("C122").Value = 2359
$ws.Range("D122").Value = 3316951
$ws.Range("C124").Value = 531313
$ws.Range("D124").Value = 701911255
$ws.Range("C126").Value = 221
$ws.Range("D126").Value = 325509
$ws.Range("C128").Value = 24
$ws.Range("D128").Value = 35998
$ws.Range("C129").Value = 1401
$ws.Range("D129").Value = 2076682
$ws.Range("C131").Value = 212945
$ws.Range("D131").Value = 313018040
$ws.Range("C132").Value = 416
$ws.Range("D132").Value = 620710
$ws.Range("C134").Value = 189756
$ws.Range("D134").Value = 275957706
$ws.Range("C137").Value = 2882
$ws.Range("D137").Value = 4046752
$ws.Range("C140").Value = 6707
$ws.Range("D140").Value = 9473603
$ws.Range("C143").Value = 45815
$ws.Range("D143").Value = 61154904
$ws.Range("C149").Value = 14345
$ws.Range("D149").Value = 21026910
$ws.Range("C150").Value = 3845
$ws.Range("D150").Value = 5543592
$ws.Range("C155").Value = 410
$ws.Range("D155").Value = 577813
$ws.Range("C156").Value = 18046
$ws.Range("D156").Value = 23858981
$ws.Range("C160").Value = 7355
$ws.Range("D160").Value = 10701429
$ws.Range("C162").Value = 5132
$ws.Range("D162").Value = 7389073
$ws.Range("C164").Value = 285
$ws.Range("D164").Value = 394139
$ws.Range("C165").Value = 277
$ws.Range("D165").Value = 395664
$ws.Range("C167").Value = 19937
$ws.Range("D167").Value = 33722794
$ws.Range("C168").Value = 2147
$ws.Range("D168").Value = 3640234
$ws.Range("C169").Value = 287
$ws.Range("D169").Value = 480089
$ws.Range("C171").Value = 67
$ws.Range("D171").Value = 119690
$ws.Range("C172").Value = 114
$ws.Range("D172").Value = 203949
$ws.Range("C173").Value = 89256
$ws.Range("D173").Value = 111525245
$ws.Range("C180").Value = 34334
$ws.Range("D180").Value = 50342558
$ws.Range("C182").Value = 13279
$ws.Range("D182").Value = 19185092
$ws.Range("C184").Value = 1258
$ws.Range("D184").Value = 1761839
$ws.Range("C186").Value = 1729
$ws.Range("D186").Value = 2426207
$ws.Range("C188").Value = 242532
$ws.Range("D188").Value = 301328956
$ws.Range("C194").Value = 888
$ws.Range("D194").Value = 1306345
$ws.Range("C196").Value = 87688
$ws.Range("D196").Value = 128522041
$ws.Range("C199").Value = 33636
$ws.Range("D199").Value = 48426042
$ws.Range("C201").Value = 17
$ws.Range("D201").Value = 23608
$ws.Range("C202").Value = 5174
$ws.Range("D202").Value = 7366396
$ws.Range("C205").Value = 5122
$ws.Range("D205").Value = 7093246
$ws.Range("C208").Value = 268856
$ws.Range("D208").Value = 332642480
$ws.Range("C209").Value = 161
$ws.Range("D209").Value = 177913
$ws.Range("C215").Value = 623
$ws.Range("D215").Value = 907378
$ws.Range("C217").Value = 96538
$ws.Range("D217").Value = 141218029
$ws.Range("C220").Value = 52572
$ws.Range("D220").Value = 75981869
$ws.Range("C223").Value = 4729
$ws.Range("D223").Value = 6637365
$ws.Range("C225").Value = 20
$ws.Range("D225").Value = 29238
$ws.Range("C226").Value = 6079
$ws.Range("D226").Value = 8420386
$ws.Range("C229").Value = 108331
$ws.Range("D229").Value = 135409601
$ws.Range("C234").Value = 568
$ws.Range("D234").Value = 829939
$ws.Range("C236").Value = 50160
$ws.Range("D236").Value = 73477831
$ws.Range("C237").Value = 38
$ws.Range("D237").Value = 54711
$ws.Range("C238").Value = 12728
$ws.Range("D238").Value = 18306818
$ws.Range("C240").Value = 1907
$ws.Range("D240").Value = 2734382
$ws.Range("C242").Value = 2634
$ws.Range("D242").Value = 3688444
$ws.Range("C243").Value = 263058
$ws.Range("D243").Value = 332122744
$ws.Range("C245").Value = 253
$ws.Range("D245").Value = 363312
$ws.Range("C248").Value = 14
$ws.Range("D248").Value = 20480
$ws.Range("C249").Value = 843
$ws.Range("D249").Value = 1237904
$ws.Range("C251").Value = 97313
$ws.Range("D251").Value = 142583196
$ws.Range("C254").Value = 66620
$ws.Range("D254").Value = 96556577
$ws.Range("C256").Value = 2447
$ws.Range("D256").Value = 3451224
$ws.Range("C259").Value = 4862
$ws.Range("D259").Value = 6828926
